$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-08-09"

# Update the column header label for the current-year total
$ws.Range("I1").Value = "2022 (through 08-09)"

# Update August's 2022 figure (row 9) and the recalculated yearly Total (row 14)
$ws.Range("I9").Value = 50
$ws.Range("I14").Value = 1020
